# Insert a new weekly record at row 46 in the "Ajo" (garlic) price sheet.
# This shifts the existing rows 46-159 down to 47-160 (preserving all of
# their data/formatting), and the new row 46 is populated with the new
# observation. The dimension (A1:R159 -> A1:R160) is updated automatically
# by Excel when the sheet is edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 46..159 down to 47..160, leaving a blank row 46 that inherits
# the formatting (notably the date-formatted style on column D) from the
# surrounding rows.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new record's data.
$ws.Cells.Item(46, 1).Value = 7
$ws.Cells.Item(46, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(46, 3).Value = "Ñuble"
$ws.Cells.Item(46, 4).Value = 44536
$ws.Cells.Item(46, 5).Value = 16
$ws.Cells.Item(46, 6).Value = 100112003
$ws.Cells.Item(46, 7).Value = "Ajo"
$ws.Cells.Item(46, 8).Value = "Chino"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 60
$ws.Cells.Item(46, 11).Value = 20000
$ws.Cells.Item(46, 12).Value = 21000
$ws.Cells.Item(46, 13).Value = 20500
$ws.Cells.Item(46, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(46, 15).Value = "China"
$ws.Cells.Item(46, 16).Value = 2050
$ws.Cells.Item(46, 17).Value = 10
$ws.Cells.Item(46, 18).Value = "Hortaliza"
